$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.52413133333333
$ws.Range("H2").Value = 37.572394
$ws.Range("I2").Value = 0.09718402715578596
$ws.Range("J2").Value = 0.1008592412859651
$ws.Range("M2").Value = 1.400501333333333
$ws.Range("N2").Value = 4.201504
$ws.Range("O2").Value = 0.00926314904242919
$ws.Range("P2").Value = 0.009687730200823723
$ws.Range("Q2").Value = 17.54006263117511
$ws.Range("R2").Value = 157.860563680576
$ws.Range("S2").Value = 0.0009002301280875312
$ws.Range("T2").Value = 0.0009770971178382113
$ws.Range("G3").Value = 12.52413133333333
$ws.Range("H3").Value = 37.572394
$ws.Range("I3").Value = 0.09718402715578596
$ws.Range("J3").Value = 0.1008592412859651
$ws.Range("O3").Value = 0.1405812059498714
$ws.Range("P3").Value = 0.1470248171880475
$ws.Range("Q3").Value = 266.1949133963446
$ws.Range("R3").Value = 2395.754220567102
$ws.Range("S3").Value = 0.01366224773662544
$ws.Range("T3").Value = 0.01482881151179419
$ws.Range("G4").Value = 12.52413133333333
$ws.Range("H4").Value = 37.572394
$ws.Range("I4").Value = 0.09718402715578596
$ws.Range("J4").Value = 0.1008592412859651
$ws.Range("M4").Value = 63.87756733333333
$ws.Range("N4").Value = 191.632702
$ws.Range("O4").Value = 0.4224968677952986
$ws.Range("P4").Value = 0.4418622271050682
$ws.Range("Q4").Value = 800.0110425365096
$ws.Range("R4").Value = 7200.099382828587
$ws.Range("S4").Value = 0.04105994707305281
$ws.Range("T4").Value = 0.044565888978744
$ws.Range("G5").Value = 12.52413133333333
$ws.Range("H5").Value = 37.572394
$ws.Range("I5").Value = 0.09718402715578596
$ws.Range("J5").Value = 0.1008592412859651
$ws.Range("M5").Value = 19.878555
$ws.Range("N5").Value = 39.75711
$ws.Range("O5").Value = 0.1314800731212866
$ws.Range("P5").Value = 0.0916710195312133
$ws.Range("Q5").Value = 248.9616335368899
$ws.Range("R5").Value = 1493.76980122134
$ws.Range("S5").Value = 0.01277776299666384
$ws.Range("T5").Value = 0.009245869477829065
$ws.Range("G6").Value = 12.52413133333333
$ws.Range("H6").Value = 37.572394
$ws.Range("I6").Value = 0.09718402715578596
$ws.Range("J6").Value = 0.1008592412859651
$ws.Range("M6").Value = 44.77944466666667
$ws.Range("N6").Value = 134.338334
$ws.Range("O6").Value = 0.2961787040911142
$ws.Range("P6").Value = 0.3097542059748472
$ws.Range("Q6").Value = 560.8236460390661
$ws.Range("R6").Value = 5047.412814351595
$ws.Range("S6").Value = 0.02878383922135634
$ws.Range("T6").Value = 0.03124157419975965
$ws.Range("I7").Value = 0.1842225641940495
$ws.Range("J7").Value = 0.1911893198517306
$ws.Range("M7").Value = 1.400501333333333
$ws.Range("N7").Value = 4.201504
$ws.Range("O7").Value = 0.00926314904242919
$ws.Range("P7").Value = 0.009687730200823723
$ws.Range("Q7").Value = 33.24903699308089
$ws.Range("R7").Value = 299.241332937728
$ws.Range("S7").Value = 0.00170648106910796
$ws.Range("T7").Value = 0.001852190548002557
$ws.Range("I8").Value = 0.1842225641940495
$ws.Range("J8").Value = 0.1911893198517306
$ws.Range("O8").Value = 0.1405812059498714
$ws.Range("P8").Value = 0.1470248171880475
$ws.Range("S8").Value = 0.02589823023757707
$ws.Range("T8").Value = 0.02810957479950783
$ws.Range("I9").Value = 0.1842225641940495
$ws.Range("J9").Value = 0.1911893198517306
$ws.Range("M9").Value = 63.87756733333333
$ws.Range("N9").Value = 191.632702
$ws.Range("O9").Value = 0.4224968677952986
$ws.Range("P9").Value = 0.4418622271050682
$ws.Range("Q9").Value = 1516.505231907918
$ws.Range("R9").Value = 13648.54708717126
$ws.Range("S9").Value = 0.07783345634920423
$ws.Range("T9").Value = 0.08447933866838893
$ws.Range("I10").Value = 0.1842225641940495
$ws.Range("J10").Value = 0.1911893198517306
$ws.Range("M10").Value = 19.878555
$ws.Range("N10").Value = 39.75711
$ws.Range("O10").Value = 0.1314800731212866
$ws.Range("P10").Value = 0.0916710195312133
$ws.Range("Q10").Value = 471.9330105819199
$ws.Range("R10").Value = 2831.598063491519
$ws.Range("S10").Value = 0.02422159621082453
$ws.Range("T10").Value = 0.01752651987428739
$ws.Range("I11").Value = 0.1842225641940495
$ws.Range("J11").Value = 0.1911893198517306
$ws.Range("M11").Value = 44.77944466666667
$ws.Range("N11").Value = 134.338334
$ws.Range("O11").Value = 0.2961787040911142
$ws.Range("P11").Value = 0.3097542059748472
$ws.Range("Q11").Value = 1063.10031758981
$ws.Range("R11").Value = 9567.902858308287
$ws.Range("S11").Value = 0.05456280032733568
$ws.Range("T11").Value = 0.0592216959615439
$ws.Range("G12").Value = 41.01852933333333
$ws.Range("H12").Value = 123.055588
$ws.Range("I12").Value = 0.3182932023406124
$ws.Range("J12").Value = 0.3303301152883236
$ws.Range("M12").Value = 1.400501333333333
$ws.Range("N12").Value = 4.201504
$ws.Range("O12").Value = 0.00926314904242919
$ws.Range("P12").Value = 0.009687730200823723
$ws.Range("Q12").Value = 57.44650502270578
$ws.Range("R12").Value = 517.018545204352
$ws.Range("S12").Value = 0.002948397372473164
$ws.Range("T12").Value = 0.003200149034120274
$ws.Range("G13").Value = 41.01852933333333
$ws.Range("H13").Value = 123.055588
$ws.Range("I13").Value = 0.3182932023406124
$ws.Range("J13").Value = 0.3303301152883236
$ws.Range("O13").Value = 0.1405812059498714
$ws.Range("P13").Value = 0.1470248171880475
$ws.Range("Q13").Value = 871.8308338456226
$ws.Range("R13").Value = 7846.477504610604
$ws.Range("S13").Value = 0.04474604223068972
$ws.Range("T13").Value = 0.04856672481197242
$ws.Range("G14").Value = 41.01852933333333
$ws.Range("H14").Value = 123.055588
$ws.Range("I14").Value = 0.3182932023406124
$ws.Range("J14").Value = 0.3303301152883236
$ws.Range("M14").Value = 63.87756733333333
$ws.Range("N14").Value = 191.632702
$ws.Range("O14").Value = 0.4224968677952986
$ws.Range("P14").Value = 0.4418622271050682
$ws.Range("Q14").Value = 2620.163869404309
$ws.Range("R14").Value = 23581.47482463878
$ws.Range("S14").Value = 0.1344778810294439
$ws.Range("T14").Value = 0.1459604004211726
$ws.Range("G15").Value = 41.01852933333333
$ws.Range("H15").Value = 123.055588
$ws.Range("I15").Value = 0.3182932023406124
$ws.Range("J15").Value = 0.3303301152883236
$ws.Range("M15").Value = 19.878555
$ws.Range("N15").Value = 39.75711
$ws.Range("O15").Value = 0.1314800731212866
$ws.Range("P15").Value = 0.0916710195312133
$ws.Range("Q15").Value = 815.38909137178
$ws.Range("R15").Value = 4892.33454823068
$ws.Range("S15").Value = 0.04184921351775218
$ws.Range("T15").Value = 0.03028169845034385
$ws.Range("G16").Value = 41.01852933333333
$ws.Range("H16").Value = 123.055588
$ws.Range("I16").Value = 0.3182932023406124
$ws.Range("J16").Value = 0.3303301152883236
$ws.Range("M16").Value = 44.77944466666667
$ws.Range("N16").Value = 134.338334
$ws.Range("O16").Value = 0.2961787040911142
$ws.Range("P16").Value = 0.3097542059748472
$ws.Range("Q16").Value = 1836.786964590044
$ws.Range("R16").Value = 16531.08268131039
$ws.Range("S16").Value = 0.09427166819025339
$ws.Range("T16").Value = 0.1023211425707144
$ws.Range("G17").Value = 14.087727
$ws.Range("H17").Value = 28.175454
$ws.Range("I17").Value = 0.1093171260259301
$ws.Range("J17").Value = 0.07563411885139956
$ws.Range("M17").Value = 1.400501333333333
$ws.Range("N17").Value = 4.201504
$ws.Range("O17").Value = 0.00926314904242919
$ws.Range("P17").Value = 0.009687730200823723
$ws.Range("Q17").Value = 19.729880447136
$ws.Range("R17").Value = 118.379282682816
$ws.Range("S17").Value = 0.001012620831268205
$ws.Range("T17").Value = 0.0007327229374093944
$ws.Range("G18").Value = 14.087727
$ws.Range("H18").Value = 28.175454
$ws.Range("I18").Value = 0.1093171260259301
$ws.Range("J18").Value = 0.07563411885139956
$ws.Range("O18").Value = 0.1405812059498714
$ws.Range("P18").Value = 0.1470248171880475
$ws.Range("Q18").Value = 299.428452872847
$ws.Range("R18").Value = 1796.570717237082
$ws.Range("S18").Value = 0.01536793340769932
$ws.Range("T18").Value = 0.01112009249730607
$ws.Range("G19").Value = 14.087727
$ws.Range("H19").Value = 28.175454
$ws.Range("I19").Value = 0.1093171260259301
$ws.Range("J19").Value = 0.07563411885139956
$ws.Range("M19").Value = 63.87756733333333
$ws.Range("N19").Value = 191.632702
$ws.Range("O19").Value = 0.4224968677952986
$ws.Range("P19").Value = 0.4418622271050682
$ws.Range("Q19").Value = 899.889730016118
$ws.Range("R19").Value = 5399.338380096709
$ws.Range("S19").Value = 0.04618614334233936
$ws.Range("T19").Value = 0.03341986020080884
$ws.Range("G20").Value = 14.087727
$ws.Range("H20").Value = 28.175454
$ws.Range("I20").Value = 0.1093171260259301
$ws.Range("J20").Value = 0.07563411885139956
$ws.Range("M20").Value = 19.878555
$ws.Range("N20").Value = 39.75711
$ws.Range("O20").Value = 0.1314800731212866
$ws.Range("P20").Value = 0.0916710195312133
$ws.Range("Q20").Value = 280.043655994485
$ws.Range("R20").Value = 1120.17462397794
$ws.Range("S20").Value = 0.01437302372329818
$ws.Range("T20").Value = 0.006933456786452757
$ws.Range("G21").Value = 14.087727
$ws.Range("H21").Value = 28.175454
$ws.Range("I21").Value = 0.1093171260259301
$ws.Range("J21").Value = 0.07563411885139956
$ws.Range("M21").Value = 44.77944466666667
$ws.Range("N21").Value = 134.338334
$ws.Range("O21").Value = 0.2961787040911142
$ws.Range("P21").Value = 0.3097542059748472
$ws.Range("Q21").Value = 630.8405916756061
$ws.Range("R21").Value = 3785.043550053636
$ws.Range("S21").Value = 0.03237740472132498
$ws.Range("T21").Value = 0.02342798642942249
$ws.Range("G22").Value = 37.49906666666667
$ws.Range("H22").Value = 112.4972
$ws.Range("I22").Value = 0.2909830802836222
$ws.Range("J22").Value = 0.3019872047225811
$ws.Range("M22").Value = 1.400501333333333
$ws.Range("N22").Value = 4.201504
$ws.Range("O22").Value = 0.00926314904242919
$ws.Range("P22").Value = 0.009687730200823723
$ws.Range("Q22").Value = 52.51749286542223
$ws.Range("R22").Value = 472.6574357888
$ws.Range("S22").Value = 0.002695419641492332
$ws.Range("T22").Value = 0.002925570563453285
$ws.Range("G23").Value = 37.49906666666667
$ws.Range("H23").Value = 112.4972
$ws.Range("I23").Value = 0.2909830802836222
$ws.Range("J23").Value = 0.3019872047225811
$ws.Range("O23").Value = 0.1405812059498714
$ws.Range("P23").Value = 0.1470248171880475
$ws.Range("Q23").Value = 797.0261999097334
$ws.Range("R23").Value = 7173.235799187601
$ws.Range("S23").Value = 0.04090675233727986
$ws.Range("T23").Value = 0.04439961356746695
$ws.Range("G24").Value = 37.49906666666667
$ws.Range("H24").Value = 112.4972
$ws.Range("I24").Value = 0.2909830802836222
$ws.Range("J24").Value = 0.3019872047225811
$ws.Range("M24").Value = 63.87756733333333
$ws.Range("N24").Value = 191.632702
$ws.Range("O24").Value = 0.4224968677952986
$ws.Range("P24").Value = 0.4418622271050682
$ws.Range("Q24").Value = 2395.349155937156
$ws.Range("R24").Value = 21558.1424034344
$ws.Range("S24").Value = 0.1229394400012583
$ws.Range("T24").Value = 0.1334367388359539
$ws.Range("G25").Value = 37.49906666666667
$ws.Range("H25").Value = 112.4972
$ws.Range("I25").Value = 0.2909830802836222
$ws.Range("J25").Value = 0.3019872047225811
$ws.Range("M25").Value = 19.878555
$ws.Range("N25").Value = 39.75711
$ws.Range("O25").Value = 0.1314800731212866
$ws.Range("P25").Value = 0.0916710195312133
$ws.Range("Q25").Value = 745.4272591820001
$ws.Range("R25").Value = 4472.563555092
$ws.Range("S25").Value = 0.03825847667274785
$ws.Range("T25").Value = 0.02768347494230024
$ws.Range("G26").Value = 37.49906666666667
$ws.Range("H26").Value = 112.4972
$ws.Range("I26").Value = 0.2909830802836222
$ws.Range("J26").Value = 0.3019872047225811
$ws.Range("M26").Value = 44.77944466666667
$ws.Range("N26").Value = 134.338334
$ws.Range("O26").Value = 0.2961787040911142
$ws.Range("P26").Value = 0.3097542059748472
$ws.Range("Q26").Value = 1679.187380851645
$ws.Range("R26").Value = 15112.6864276648
$ws.Range("S26").Value = 0.08618299163084389
$ws.Range("T26").Value = 0.09354180681340672
